# Slide 8: "Phoenix, AZ / Employment and Home Prices" -> retitle, re-layout the
# title/body/picture, and replace the two short bullet paragraphs with one
# longer written-out paragraph.
#
# NOTE on the Left/Top/Width/Height numbers below: PowerPoint's Shape
# position/size properties are expressed in points (Single-precision) while
# the underlying OOXML stores EMU (1 pt = 12700 EMU). The literals here are
# chosen so that pt -> EMU round-trips to the exact target EMU offsets/
# extents from the source file once run through that point->EMU conversion.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- Title placeholder: reposition, enable shrink-text-on-overflow, update text ---
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Runs(2).Text = "Employment and Median Home Sale Prices"
$title.TextFrame.AutoSize = 2
$title.Left = 113.50008014015748
$title.Top = 17.095511811023623
$title.Width = 756.0
$title.Height = 117.0

# --- Body textbox: new combined paragraph, font size 16pt, reposition/resize ---
$tb = $s.Shapes.Item(2)
$tb.TextFrame.TextRange.Text = "Before 2009, the median sales price was about 225,000 which is almost hit the highest point in the recently 10 years. The employment number in phoenix area before 2009, which was before the Great Recession, was about 1900 thousands still not bad. Between 2009 to 2013, the Great Recession was happening, more and more people lost the job and they can’t afford the loan of the houses. There were lots of foreclosure houses in the market, the price of houses was dropped down to the lowest point in the recently 10 years, and the home prices even fell faster than employment. After 2013, which after the Great Recession, there were more and more people find the job, employment number was raising up, home prices are also raising up."
$tb.TextFrame.TextRange.Font.Size = 16
$tb.Left = 110.60559055118111
$tb.Top = 349.5677952755905
$tb.Width = 735.5
$tb.Height = 162.3703149606299

# --- Picture: reposition/resize ---
$pic = $s.Shapes.Item(3)
$pic.Left = 110.60559055118111
$pic.Top = 124.05464566929133
$pic.Width = 534.1399231598425
$pic.Height = 222.5583464566929
